# Insert a new weekly price record for "Haba" (Femacal de La Calera) above
# the existing row 255, shifting the rest of the table (old rows 255-289)
# down by one to rows 256-290.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 255; all rows below shift down.
$ws.Rows(255).Insert()

# Populate the newly inserted row 255 with the new record.
$ws.Cells.Item(255, 1).Value = 3
$ws.Cells.Item(255, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(255, 3).Value = "Coquimbo"
$ws.Cells.Item(255, 4).Value = 45154
$ws.Cells.Item(255, 5).Value = 5
$ws.Cells.Item(255, 6).Value = 100112026
$ws.Cells.Item(255, 7).Value = "Haba"
$ws.Cells.Item(255, 8).Value = "Sin especificar"
$ws.Cells.Item(255, 9).Value = "Primera"
$ws.Cells.Item(255, 10).Value = 45
$ws.Cells.Item(255, 11).Value = 16000
$ws.Cells.Item(255, 12).Value = 16000
$ws.Cells.Item(255, 13).Value = 16000
$ws.Cells.Item(255, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(255, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(255, 16).Value = 640
$ws.Cells.Item(255, 17).Value = 25
$ws.Cells.Item(255, 18).Value = "Hortaliza"
